# Auto-generated Excel COM-interop script
# Applies the "Phantom_Profits" numeric corrections to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 361.7143  # H33
$ws.Cells.Item(33, 9).Value = 466.4  # I33
$ws.Cells.Item(33, 10).Value = 100  # J33
$ws.Cells.Item(33, 11).Value = 466.4  # K33
$ws.Cells.Item(33, 12).Value = 100  # L33
$ws.Cells.Item(33, 13).Value = -237.4  # M33
$ws.Cells.Item(33, 14).Value = -558  # N33
$ws.Cells.Item(51, 8).Value = 10572.091  # H51
$ws.Cells.Item(51, 10).Value = 11399.6  # J51
$ws.Cells.Item(51, 12).Value = 11399.6  # L51
$ws.Cells.Item(51, 14).Value = -12367.6  # N51
$ws.Cells.Item(69, 8).Value = 17317.072  # H69
$ws.Cells.Item(69, 10).Value = 18187.615  # J69
$ws.Cells.Item(69, 12).Value = 54562.845  # L69
$ws.Cells.Item(69, 14).Value = -56310.845  # N69
$ws.Cells.Item(70, 8).Value = 10157.667  # H70
$ws.Cells.Item(70, 10).Value = 9989.200000000001  # J70
$ws.Cells.Item(70, 12).Value = 29967.6  # L70
$ws.Cells.Item(70, 14).Value = -30507.6  # N70
$ws.Cells.Item(72, 8).Value = 17317.072  # H72
$ws.Cells.Item(72, 10).Value = 18187.615  # J72
$ws.Cells.Item(72, 12).Value = 163688.535  # L72
$ws.Cells.Item(72, 14).Value = -172424.535  # N72
$ws.Cells.Item(73, 8).Value = 10157.667  # H73
$ws.Cells.Item(73, 10).Value = 9989.200000000001  # J73
$ws.Cells.Item(73, 12).Value = 29967.6  # L73
$ws.Cells.Item(73, 14).Value = -31839.6  # N73
$ws.Cells.Item(98, 8).Value = 1403  # H98
$ws.Cells.Item(98, 9).Value = 578.375  # I98
$ws.Cells.Item(98, 11).Value = 578.375  # K98
$ws.Cells.Item(98, 13).Value = 919.625  # M98
$ws.Cells.Item(122, 8).Value = 1403  # H122
$ws.Cells.Item(122, 9).Value = 578.375  # I122
$ws.Cells.Item(122, 11).Value = 1735.125  # K122
$ws.Cells.Item(122, 13).Value = 714.875  # M122
$ws.Cells.Item(132, 8).Value = 813  # H132
$ws.Cells.Item(132, 9).Value = 841.5909  # I132
$ws.Cells.Item(132, 11).Value = 2524.7727  # K132
$ws.Cells.Item(132, 13).Value = 5.227300000000014  # M132
$ws.Cells.Item(138, 8).Value = 2299.3333  # H138
$ws.Cells.Item(138, 9).Value = 1377.5555  # I138
$ws.Cells.Item(138, 11).Value = 4132.666499999999  # K138
$ws.Cells.Item(138, 13).Value = 1007.333500000001  # M138
$ws.Cells.Item(116, 8).Value = 5000  # H116
$ws.Cells.Item(116, 9).Value = 5000  # I116
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 11).Value = 5000  # K116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 13).Value = -1558  # M116
$ws.Cells.Item(116, 14).ClearContents()  # N116

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4619.9688  # H32
$ws.Cells.Item(32, 9).Value = 4619.484  # I32
$ws.Cells.Item(32, 11).Value = 4619.484  # K32
$ws.Cells.Item(32, 13).Value = -4332.484  # M32
$ws.Cells.Item(61, 8).Value = 3118.15  # H61
$ws.Cells.Item(61, 9).Value = 2920.1667  # I61
$ws.Cells.Item(61, 11).Value = 2920.1667  # K61
$ws.Cells.Item(61, 13).Value = -2708.1667  # M61
$ws.Cells.Item(97, 8).Value = 1833.375  # H97
$ws.Cells.Item(97, 9).Value = 666.7143  # I97
$ws.Cells.Item(97, 11).Value = 666.7143  # K97
$ws.Cells.Item(97, 13).Value = -170.7143  # M97
$ws.Cells.Item(101, 8).Value = 31583.166  # H101
$ws.Cells.Item(101, 10).Value = 31583.166  # J101
$ws.Cells.Item(101, 12).Value = 31583.166  # L101
$ws.Cells.Item(101, 14).Value = -38073.166  # N101
$ws.Cells.Item(136, 8).Value = 3118.15  # H136
$ws.Cells.Item(136, 9).Value = 2920.1667  # I136
$ws.Cells.Item(136, 11).Value = 8760.500100000001  # K136
$ws.Cells.Item(136, 13).Value = -6210.500100000001  # M136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1656.3125  # H20
$ws.Cells.Item(20, 9).Value = 1527.909  # I20
$ws.Cells.Item(20, 11).Value = 1527.909  # K20
$ws.Cells.Item(20, 13).Value = -1280.909  # M20
$ws.Cells.Item(134, 8).Value = 5978.2666  # H134
$ws.Cells.Item(134, 9).Value = 5697.4614  # I134
$ws.Cells.Item(134, 11).Value = 17092.3842  # K134
$ws.Cells.Item(134, 13).Value = -14557.3842  # M134
$ws.Cells.Item(22, 8).Value = 649  # H22
$ws.Cells.Item(22, 9).Value = 465.33334  # I22
$ws.Cells.Item(22, 10).Value = 1200  # J22
$ws.Cells.Item(22, 11).Value = 465.33334  # K22
$ws.Cells.Item(22, 12).Value = 1200  # L22
$ws.Cells.Item(22, 13).Value = -292.33334  # M22
$ws.Cells.Item(22, 14).Value = -1546  # N22
$ws.Cells.Item(105, 8).Value = 4496.1665  # H105
$ws.Cells.Item(105, 9).Value = 4496.1665  # I105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 11).Value = 4496.1665  # K105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 13).Value = -2749.1665  # M105
$ws.Cells.Item(105, 14).ClearContents()  # N105

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2225.8696  # H58
$ws.Cells.Item(58, 9).Value = 1795.0476  # I58
$ws.Cells.Item(58, 11).Value = 1795.0476  # K58
$ws.Cells.Item(58, 13).Value = -1592.0476  # M58
$ws.Cells.Item(86, 8).Value = 9996.75  # H86
$ws.Cells.Item(86, 9).Value = 9995.666999999999  # I86
$ws.Cells.Item(86, 10).Value = 10000  # J86
$ws.Cells.Item(86, 11).Value = 9995.666999999999  # K86
$ws.Cells.Item(86, 12).Value = 10000  # L86
$ws.Cells.Item(86, 13).Value = -8872.666999999999  # M86
$ws.Cells.Item(86, 14).Value = -12246  # N86
$ws.Cells.Item(89, 8).Value = 9996.75  # H89
$ws.Cells.Item(89, 9).Value = 9995.666999999999  # I89
$ws.Cells.Item(89, 10).Value = 10000  # J89
$ws.Cells.Item(89, 11).Value = 49978.335  # K89
$ws.Cells.Item(89, 12).Value = 50000  # L89
$ws.Cells.Item(89, 13).Value = -44362.335  # M89
$ws.Cells.Item(89, 14).Value = -61232  # N89
$ws.Cells.Item(94, 8).Value = 1526.3334  # H94
$ws.Cells.Item(94, 10).Value = 1999  # J94
$ws.Cells.Item(94, 12).Value = 1999  # L94
$ws.Cells.Item(94, 14).Value = -2901  # N94
$ws.Cells.Item(99, 8).Value = 2231.625  # H99
$ws.Cells.Item(99, 9).Value = 2409  # I99
$ws.Cells.Item(99, 11).Value = 2409  # K99
$ws.Cells.Item(99, 13).Value = -911  # M99
$ws.Cells.Item(126, 8).Value = 2231.625  # H126
$ws.Cells.Item(126, 9).Value = 2409  # I126
$ws.Cells.Item(126, 11).Value = 7227  # K126
$ws.Cells.Item(126, 13).Value = -4757  # M126
$ws.Cells.Item(136, 8).Value = 2225.8696  # H136
$ws.Cells.Item(136, 9).Value = 1795.0476  # I136
$ws.Cells.Item(136, 11).Value = 5385.142800000001  # K136
$ws.Cells.Item(136, 13).Value = -2835.142800000001  # M136

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 57.307693  # H12
$ws.Cells.Item(12, 9).Value = 56.8  # I12
$ws.Cells.Item(12, 10).Value = 57.625  # J12
$ws.Cells.Item(12, 11).Value = 170.4  # K12
$ws.Cells.Item(12, 12).Value = 172.875  # L12
$ws.Cells.Item(12, 13).Value = 2.600000000000023  # M12
$ws.Cells.Item(12, 14).Value = -518.875  # N12
$ws.Cells.Item(62, 8).Value = 10725  # H62
$ws.Cells.Item(62, 10).Value = 10787.5  # J62
$ws.Cells.Item(62, 12).Value = 32362.5  # L62
$ws.Cells.Item(62, 14).Value = -33734.5  # N62
$ws.Cells.Item(65, 8).Value = 10725  # H65
$ws.Cells.Item(65, 10).Value = 10787.5  # J65
$ws.Cells.Item(65, 12).Value = 97087.5  # L65
$ws.Cells.Item(65, 14).Value = -103951.5  # N65
$ws.Cells.Item(140, 8).Value = 716699.4399999999  # H140
$ws.Cells.Item(140, 9).Value = 716699.4399999999  # I140
$ws.Cells.Item(140, 11).Value = 2150098.32  # K140
$ws.Cells.Item(140, 13).Value = -2144918.32  # M140
$ws.Cells.Item(82, 8).Value = 1200  # H82
$ws.Cells.Item(82, 9).Value = 1200  # I82
$ws.Cells.Item(82, 10).Value = 0  # J82
$ws.Cells.Item(82, 11).Value = 3600  # K82
$ws.Cells.Item(82, 12).Value = 0  # L82
$ws.Cells.Item(82, 14).ClearContents()  # N82
$ws.Cells.Item(82, 13).Value = -3194  # M82
$ws.Cells.Item(85, 8).Value = 1200  # H85
$ws.Cells.Item(85, 9).Value = 1200  # I85
$ws.Cells.Item(85, 10).Value = 0  # J85
$ws.Cells.Item(85, 11).Value = 3600  # K85
$ws.Cells.Item(85, 12).Value = 0  # L85
$ws.Cells.Item(85, 14).ClearContents()  # N85
$ws.Cells.Item(85, 13).Value = -2196  # M85

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 126.5  # H2
$ws.Cells.Item(2, 9).Value = 133.26666  # I2
$ws.Cells.Item(2, 11).Value = 133.26666  # K2
$ws.Cells.Item(2, 13).Value = -20.26666  # M2
$ws.Cells.Item(97, 8).Value = 605.8  # H97
$ws.Cells.Item(97, 9).Value = 550.44446  # I97
$ws.Cells.Item(97, 11).Value = 550.44446  # K97
$ws.Cells.Item(97, 13).Value = -54.44446000000005  # M97
$ws.Cells.Item(102, 8).Value = 1323  # H102
$ws.Cells.Item(102, 9).Value = 1323  # I102
$ws.Cells.Item(102, 11).Value = 1323  # K102
$ws.Cells.Item(102, 13).Value = 299  # M102

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1490.7142  # H22
$ws.Cells.Item(22, 9).Value = 1416  # I22
$ws.Cells.Item(22, 11).Value = 1416  # K22
$ws.Cells.Item(22, 13).Value = -1121  # M22
$ws.Cells.Item(27, 8).Value = 1490.7142  # H27
$ws.Cells.Item(27, 9).Value = 1416  # I27
$ws.Cells.Item(27, 11).Value = 1416  # K27
$ws.Cells.Item(27, 13).Value = -1309  # M27
$ws.Cells.Item(122, 8).Value = 3407.2727  # H122
$ws.Cells.Item(122, 9).Value = 3348.4  # I122
$ws.Cells.Item(122, 11).Value = 10045.2  # K122
$ws.Cells.Item(122, 13).Value = -7595.200000000001  # M122
$ws.Cells.Item(93, 8).Value = 2999  # H93
$ws.Cells.Item(93, 10).Value = 2998  # J93
$ws.Cells.Item(93, 12).Value = 2998  # L93
$ws.Cells.Item(93, 14).Value = -5494  # N93

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 18111.875  # H74
$ws.Cells.Item(74, 9).Value = 17732  # I74
$ws.Cells.Item(74, 11).Value = 17732  # K74
$ws.Cells.Item(74, 13).Value = -16796  # M74
$ws.Cells.Item(77, 8).Value = 18111.875  # H77
$ws.Cells.Item(77, 9).Value = 17732  # I77
$ws.Cells.Item(77, 11).Value = 53196  # K77
$ws.Cells.Item(77, 13).Value = -48516  # M77
$ws.Cells.Item(132, 8).Value = 5325.269  # H132
$ws.Cells.Item(132, 9).Value = 4252.304  # I132
$ws.Cells.Item(132, 10).Value = 13551.333  # J132
$ws.Cells.Item(132, 11).Value = 12756.912  # K132
$ws.Cells.Item(132, 12).Value = 40653.999  # L132
$ws.Cells.Item(132, 13).Value = -10226.912  # M132
$ws.Cells.Item(132, 14).Value = -45713.999  # N132
$ws.Cells.Item(136, 8).Value = 12880.4  # H136
$ws.Cells.Item(136, 9).Value = 13567.167  # I136
$ws.Cells.Item(136, 11).Value = 40701.501  # K136
$ws.Cells.Item(136, 13).Value = -38151.501  # M136

